$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.295.23"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.98"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.34"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4647"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3917"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07888"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9895"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.920.31"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.092"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06991"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.43"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001002"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.298.39"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.323"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.096"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.59"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.986"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.74"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.920"
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09375"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9078"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.292"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.226"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.182"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05795"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9995"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.785"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5727"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1788"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.787"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.00"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.206"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5357"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.595"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.33"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.34"
$ws.Range("E51").Value = "  -0.46%  "
